$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G6").Value = 2.1
$ws.Range("I6").Value = 3.75
$ws.Range("M6").Value = 1.11
$ws.Range("N6").Value = 6.5
$ws.Range("O6").Value = 1.44
$ws.Range("P6").Value = 2.63
$ws.Range("X6").Value = 1.17
$ws.Range("AI6").Value = 6.5
$ws.Range("AQ6").Value = 34
$ws.Range("M7").Value = 1.08
$ws.Range("O7").Value = 1.44
$ws.Range("P7").Value = 2.63
$ws.Range("X7").Value = 1.2
$ws.Range("G10").Value = 1.65
$ws.Range("H10").Value = 3.1
$ws.Range("I10").Value = 6.4
$ws.Range("J10").Value = 2.32
$ws.Range("K10").Value = 1.87
$ws.Range("L10").Value = 6.8
$ws.Range("M10").Value = 1.15
$ws.Range("N10").Value = 4.85
$ws.Range("O10").Value = 1.62
$ws.Range("P10").Value = 2.15
$ws.Range("S10").Value = 2.82
$ws.Range("T10").Value = 1.38
$ws.Range("W10").Value = 5.2
$ws.Range("X10").Value = 1.13
$ws.Range("Y10").Value = 1.62
$ws.Range("Z10").Value = 2.15
$ws.Range("AA10").Value = 2.65
$ws.Range("AB10").Value = 1.42
$ws.Range("AC10").Value = 4.2
$ws.Range("AD10").Value = 6
$ws.Range("AF10").Value = 12
$ws.Range("AG10").Value = 19.5
$ws.Range("AH10").Value = 60
$ws.Range("AI10").Value = 4.85
$ws.Range("AJ10").Value = 6.7
$ws.Range("AK10").Value = 28
$ws.Range("AL10").Value = 250
$ws.Range("AM10").Value = 10.5
$ws.Range("AN10").Value = 37
$ws.Range("AO10").Value = 24
$ws.Range("M11").Value = 1.08
$ws.Range("O11").Value = 1.5
$ws.Range("P11").Value = 2.63
$ws.Range("X11").Value = 1.2
$ws.Range("Y11").Value = 1.53
$ws.Range("Z11").Value = 2.38
$ws.Range("M12").Value = 1.06
$ws.Range("O12").Value = 1.33
$ws.Range("P12").Value = 3.4
$ws.Range("S12").Value = 2.05
$ws.Range("T12").Value = 1.8
$ws.Range("X12").Value = 1.3
$ws.Range("Y12").Value = 1.44
$ws.Range("Z12").Value = 2.63
$ws.Range("M13").Value = 1.1
$ws.Range("O13").Value = 1.53
$ws.Range("P13").Value = 2.5
$ws.Range("Q13").Value = 2
$ws.Range("R13").Value = 1.85
$ws.Range("X13").Value = 1.17
$ws.Range("Y13").Value = 1.62
$ws.Range("M14").Value = 1.08
$ws.Range("O14").Value = 1.44
$ws.Range("P14").Value = 2.75
$ws.Range("X14").Value = 1.22
$ws.Range("Y14").Value = 1.5
$ws.Range("G15").Value = 1.86
$ws.Range("M15").Value = 1.03
$ws.Range("O15").Value = 1.18
$ws.Range("X15").Value = 1.5
$ws.Range("Y15").Value = 1.3
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = 3.3
$ws.Range("I16").Value = 1.86
$ws.Range("J16").Value = 4.5
$ws.Range("L16").Value = 2.63
$ws.Range("AC16").Value = 11
$ws.Range("AD16").Value = 21
$ws.Range("G17").Value = 2.82
$ws.Range("I17").Value = 2.4
$ws.Range("AC17").Value = 7
$ws.Range("AJ17").Value = 6.5
$ws.Range("G18").Value = 1.71
$ws.Range("I18").Value = 4.75
$ws.Range("J18").Value = 2.4
$ws.Range("L18").Value = 4.75
$ws.Range("M18").Value = 1.05
$ws.Range("N18").Value = 11
$ws.Range("O18").Value = 1.29
$ws.Range("P18").Value = 3.5
$ws.Range("AI18").Value = 11
$ws.Range("AJ18").Value = 7
$ws.Range("AP18").Value = 51
$ws.Range("AS18").Value = 251
$ws.Range("G19").Value = 2.45
$ws.Range("I19").Value = 2.75
$ws.Range("N19").Value = 8.5
$ws.Range("Y19").Value = 1.5
$ws.Range("Z19").Value = 2.5
$ws.Range("AM19").Value = 8
$ws.Range("G20").Value = 1.66
$ws.Range("M20").Value = 1.05
$ws.Range("N20").Value = 11
$ws.Range("N22").Value = 9
$ws.Range("S22").Value = 2.1
$ws.Range("T22").Value = 1.7
$ws.Range("W22").Value = 3.75
$ws.Range("X22").Value = 1.25
$ws.Range("M23").Value = 1.06
$ws.Range("N23").Value = 10
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = 3.25
$ws.Range("I24").Value = 3.7
$ws.Range("J24").Value = 2.75
$ws.Range("L24").Value = 4.33
$ws.Range("N24").Value = 8.5
$ws.Range("AJ24").Value = 6
$ws.Range("S25").Value = 1.73
$ws.Range("T25").Value = 2.08
$ws.Range("M26").Value = 1.06
$ws.Range("N26").Value = 10
$ws.Range("S26").Value = 2.05
$ws.Range("T26").Value = 1.75
$ws.Range("I31").Value = 1.45
$ws.Range("S31").Value = 1.65
$ws.Range("G32").Value = 3.7
$ws.Range("S32").Value = 1.8
$ws.Range("T32").Value = 2
$ws.Range("W32").Value = 2.75
$ws.Range("X32").Value = 1.4
$ws.Range("AC32").Value = 13
$ws.Range("AR32").Value = 23
$ws.Range("O33").Value = 1.33
$ws.Range("P33").Value = 3.25
$ws.Range("S33").Value = 2.05
$ws.Range("AC33").Value = 9
$ws.Range("AP10").Value = 175
$ws.Range("AQ10").Value = 120
$ws.Range("AR10").Value = 150
